$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the "Rules String Hello (Integer hour)" decision table (B3:E11)
# into a second copy at B27:E35 (values, styles, number formats and the
# B3:E3-style merged header all come along with the copy).
$ws.Range("B3:E11").Copy($ws.Range("B27"))

# Give the new table its own header text (creates a new shared string).
$ws.Range("B27").Value = "Rules String Hello2 (Integer hour)"

# Rows that have a thick bottom border in the source table (3,4,6,7,11)
# also carry an explicit 13.5pt row height; reproduce that on their
# counterparts in the new block (27,28,30,31,35).
$ws.Rows.Item(27).RowHeight = 13.5
$ws.Rows.Item(28).RowHeight = 13.5
$ws.Rows.Item(30).RowHeight = 13.5
$ws.Rows.Item(31).RowHeight = 13.5
$ws.Rows.Item(35).RowHeight = 13.5

# Re-create the explanatory cell comments for the new header/condition/
# return columns (comments aren't brought over by Copy/PasteSpecial).
$ws.Range("B27").AddComment("This is so-called Decision Table Header. It starts with the keyword ""Rules"".")
$ws.Range("B28").AddComment("`nRule column header. Rule column is used to to name particular rule rows for documentation and tracing purposes. It is also useful to create rule rows that span more than one cell vertically (this will be explained in one of the next tutorials)`n")
$ws.Range("C28").AddComment("Condition column header. Must start with ""C""")
$ws.Range("E28").AddComment("Return column header. Must start with ""RET"".  ")
$ws.Range("C29").AddComment("Condition expression. Must have type boolean. As you can see condition uses parameter hour from Method Header and variable min that defines column data. When condition is evaluated for each row, the cell value from this row is assigned to variable min")
$ws.Range("E29").AddComment("This is return expression performed for the first row where all conditions have been satisfied. The variable greeting is substittuted with a cell value from the rule row")

# Match the final selection left behind in the edited file.
$ws.Range("F24").Select()
